$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

# ------------------------------------------------------------------
# 1) The login test case (TC_001, row 2) status flips from "Yes" to
#    "No". In this workbook the SmokeTest/RegressionTest columns for
#    every test case (rows 2-7) share the very same string value, so
#    changing that shared text turns all of D2:E7 into "No" as well.
# ------------------------------------------------------------------
$ws.Range("D2:E7").Value = "No"

# ------------------------------------------------------------------
# 2) Rows 4-7 (D4:E7) were using a slightly different cell style than
#    rows 2-3 (D2:E3). Bring their formatting in line by copying the
#    formatting already used on D2 onto D4:E7 (values are untouched).
# ------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("D4:E7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3) Move the active selection to E10.
# ------------------------------------------------------------------
[void]$ws.Range("E10").Select()

# ------------------------------------------------------------------
# 4) Extend the "Yes,No" list-validation so it also covers D1 and
#    D4:E7 (previously only D1:D3, E2:E3, D8:E9 and D13:E21 had it).
#    Rebuild it as a single rule across the whole D1:E21 block and
#    then carve back out the cells that must keep their own rule
#    (E1 keeps the "Yes"-only validation; D10:E12 keep no validation).
# ------------------------------------------------------------------
$ws.Range("D1:D3").Validation.Delete()
$ws.Range("E2:E3").Validation.Delete()
$ws.Range("D8:E9").Validation.Delete()
$ws.Range("D13:E21").Validation.Delete()

$ws.Range("D1:E21").Validation.Add(3, 1, 1, '"Yes,No"')

$ws.Range("E1").Validation.Delete()
$ws.Range("D10:D12").Validation.Delete()
$ws.Range("E10:E12").Validation.Delete()

$ws.Range("E1").Validation.Add(3, 1, 1, '"Yes"')
